$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 797, shifting the existing rows 797-838 down to 798-839
$ws.Rows.Item(797).Insert()

# Populate the newly inserted row with the new daily record.
# Force the date column to text format first so "2026/02/15" is stored as a
# literal string (matching the rest of the column) instead of being
# auto-parsed into a date serial number.
$ws.Cells.Item(797, 1).NumberFormat = "@"
$ws.Cells.Item(797, 1).Value = "2026/02/15"
$ws.Cells.Item(797, 1).Style = "Normal"
$ws.Cells.Item(797, 2).Value = "日"
$ws.Cells.Item(797, 3).Value = 8
$ws.Cells.Item(797, 4).Value = 23
